$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet ("Address" -> "Hoja 1")
$ws.Name = "Hoja 1"

# Fix the text of A1 (drop the German "ß" for "ss") and B2 (drop the
# trailing line break after "..., Berlin")
$ws.Range("A1").Value = "Bayerische Strasse 25"
$ws.Range("B2").Value = "Düsseldorfer Straße 3, Berlin"

# Resize the two data columns
$ws.Columns.Item(1).ColumnWidth = 19.85
$ws.Columns.Item(2).ColumnWidth = 27.5

# Row 2 no longer needs its custom height - AutoFit drops back to default
$ws.Rows.Item(2).AutoFit()

# Drop the now-empty row 5 entirely
$ws.Rows.Item(5).Delete()
